$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "real" burndown series (row 4): correct values through 03 Feb, extend two more days ---
$ws.Range("O4").Value = 73
$ws.Range("P4").Value = 70
$ws.Range("Q4").Value = 64
$ws.Range("R4").Value = 46

# --- "finished" (actual hours) series (row 6): correct 31 Jan, extend two more days ---
$ws.Range("N6").Value = 9
$ws.Range("P6").Value = 6
$ws.Range("Q6").Value = 18

# --- remove the stale per-person "week 2" (column D) breakdown, rows 11-20 (incl. the SUM formula) ---
$ws.Range("D11:D20").ClearContents()

# --- drop the trailing "TODO:" note rows ---
$ws.Rows("27:28").Delete()

# --- restore the cursor position recorded in the saved file ---
$ws.Range("D26").Select()
